$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 3000
$ws.Range("M31").Value = -2770

$ws.Range("H64").Value = 2720
$ws.Range("I64").Value = 2675
$ws.Range("K64").Value = 2675
$ws.Range("M64").Value = -2427

$ws.Range("H67").Value = 2720
$ws.Range("I67").Value = 2675
$ws.Range("K67").Value = 2675
$ws.Range("M67").Value = -1817

$ws.Range("H86").Value = 17365.13
$ws.Range("I86").Value = 9107.538
$ws.Range("K86").Value = 9107.538
$ws.Range("M86").Value = -7984.538

$ws.Range("H89").Value = 17365.13
$ws.Range("I89").Value = 9107.538
$ws.Range("K89").Value = 45537.69
$ws.Range("M89").Value = -39921.69

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = ""

$ws.Range("H121").Value = 747.7
$ws.Range("J121").Value = 922.1429000000001
$ws.Range("L121").Value = 2766.4287
$ws.Range("N121").Value = -6260.4287

$ws.Range("H132").Value = 6852783.5
$ws.Range("I132").Value = 11366762
$ws.Range("J132").Value = 3989.276
$ws.Range("K132").Value = 34100286
$ws.Range("L132").Value = 11967.828
$ws.Range("M132").Value = -34097756
$ws.Range("N132").Value = -17027.828

$ws.Range("H138").Value = 1845.5079
$ws.Range("I138").Value = 863.9729599999999
$ws.Range("J138").Value = 3242.3076
$ws.Range("K138").Value = 2591.91888
$ws.Range("L138").Value = 9726.9228
$ws.Range("M138").Value = 2548.08112
$ws.Range("N138").Value = -20006.9228

$ws.Range("H141").Value = 2373.3125
$ws.Range("I141").Value = 1078.8846
$ws.Range("J141").Value = 7982.5
$ws.Range("K141").Value = 3236.6538
$ws.Range("L141").Value = 23947.5
$ws.Range("M141").Value = 1943.3462
$ws.Range("N141").Value = -34307.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10002710
$ws.Range("I32").Value = 2764.518
$ws.Range("J32").Value = 58825972
$ws.Range("K32").Value = 2764.518
$ws.Range("L32").Value = 58825972
$ws.Range("M32").Value = -2477.518
$ws.Range("N32").Value = -58826546

$ws.Range("H61").Value = 9261166
$ws.Range("I61").Value = 10418562
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 10418562
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -10418350
$ws.Range("N61").Value = -2424

$ws.Range("H88").Value = 2799.2856
$ws.Range("I88").Value = 2759
$ws.Range("J88").Value = 2900
$ws.Range("K88").Value = 2759
$ws.Range("L88").Value = 2900
$ws.Range("M88").Value = -2353
$ws.Range("N88").Value = -3712

$ws.Range("H91").Value = 2799.2856
$ws.Range("I91").Value = 2759
$ws.Range("J91").Value = 2900
$ws.Range("K91").Value = 2759
$ws.Range("L91").Value = 2900
$ws.Range("M91").Value = -1355
$ws.Range("N91").Value = -5708

$ws.Range("H102").Value = 1112.909
$ws.Range("I102").Value = 1046.6666
$ws.Range("J102").Value = 1411
$ws.Range("K102").Value = 1046.6666
$ws.Range("L102").Value = 1411
$ws.Range("M102").Value = 575.3334
$ws.Range("N102").Value = -4655

$ws.Range("H136").Value = 9261166
$ws.Range("I136").Value = 10418562
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 31255686
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -31253136
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 896774.3
$ws.Range("I86").Value = 2690.5
$ws.Range("J86").Value = 1663131.9
$ws.Range("K86").Value = 2690.5
$ws.Range("L86").Value = 1663131.9
$ws.Range("M86").Value = -1567.5
$ws.Range("N86").Value = -1665377.9

$ws.Range("H89").Value = 896774.3
$ws.Range("I89").Value = 2690.5
$ws.Range("J89").Value = 1663131.9
$ws.Range("K89").Value = 13452.5
$ws.Range("L89").Value = 8315659.5
$ws.Range("M89").Value = -7836.5
$ws.Range("N89").Value = -8326891.5

$ws.Range("H107").Value = 23810352
$ws.Range("I107").Value = 27778630
$ws.Range("J107").Value = 678.6667
$ws.Range("K107").Value = 27778630
$ws.Range("L107").Value = 678.6667
$ws.Range("M107").Value = -27776710
$ws.Range("N107").Value = -4518.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1333.125
$ws.Range("I35").Value = 837.8570999999999
$ws.Range("J35").Value = 4800
$ws.Range("K35").Value = 837.8570999999999
$ws.Range("L35").Value = 4800
$ws.Range("M35").Value = -543.8570999999999
$ws.Range("N35").Value = -5388

$ws.Range("H62").Value = 3724.4546
$ws.Range("I62").Value = 2519.8
$ws.Range("J62").Value = 4728.3335
$ws.Range("K62").Value = 2519.8
$ws.Range("L62").Value = 4728.3335
$ws.Range("M62").Value = -1895.8
$ws.Range("N62").Value = -5976.3335

$ws.Range("H65").Value = 3724.4546
$ws.Range("I65").Value = 2519.8
$ws.Range("J65").Value = 4728.3335
$ws.Range("K65").Value = 12599
$ws.Range("L65").Value = 23641.6675
$ws.Range("M65").Value = -9479
$ws.Range("N65").Value = -29881.6675

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1153
$ws.Range("I57").Value = 1000
$ws.Range("J57").Value = 1306
$ws.Range("K57").Value = 3000
$ws.Range("L57").Value = 3918
$ws.Range("M57").Value = -2441
$ws.Range("N57").Value = -5036

$ws.Range("H120").Value = 13425.556
$ws.Range("I120").Value = 830
$ws.Range("J120").Value = 15000
$ws.Range("K120").Value = 2490
$ws.Range("L120").Value = 45000
$ws.Range("M120").Value = 2348
$ws.Range("N120").Value = -54676

$ws.Range("H131").Value = 872.61
$ws.Range("J131").Value = 881.2347
$ws.Range("L131").Value = 2643.7041
$ws.Range("N131").Value = -12723.7041

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7695785
$ws.Range("I80").Value = 5116.6665
$ws.Range("J80").Value = 14287786
$ws.Range("K80").Value = 5116.6665
$ws.Range("L80").Value = 14287786
$ws.Range("M80").Value = -4118.6665
$ws.Range("N80").Value = -14289782

$ws.Range("H83").Value = 7695785
$ws.Range("I83").Value = 5116.6665
$ws.Range("J83").Value = 14287786
$ws.Range("K83").Value = 25583.3325
$ws.Range("L83").Value = 71438930
$ws.Range("M83").Value = -20591.3325
$ws.Range("N83").Value = -71448914

$ws.Range("H132").Value = 2382.1538
$ws.Range("I132").Value = 1914
$ws.Range("J132").Value = 3266.4443
$ws.Range("K132").Value = 5742
$ws.Range("L132").Value = 9799.332900000001
$ws.Range("M132").Value = -3212
$ws.Range("N132").Value = -14859.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 11450.4
$ws.Range("I122").Value = 12556
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 37668
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -35218
$ws.Range("N122").Value = -9400

$ws.Range("H132").Value = 7069.1836
$ws.Range("I132").Value = 2229.4583
$ws.Range("J132").Value = 11715.32
$ws.Range("K132").Value = 6688.374899999999
$ws.Range("L132").Value = 35145.96
$ws.Range("M132").Value = -4158.374899999999
$ws.Range("N132").Value = -40205.96

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 4361
$ws.Range("J29").Value = 4361
$ws.Range("L29").Value = 4361
$ws.Range("N29").Value = -4941

$ws.Range("H132").Value = 20759.963
$ws.Range("I132").Value = 28260.256
$ws.Range("J132").Value = 2478
$ws.Range("K132").Value = 84780.76800000001
$ws.Range("L132").Value = 7434
$ws.Range("M132").Value = -82250.76800000001
$ws.Range("N132").Value = -12494
